$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6583.25
$ws.Range("I40").Value = 3333
$ws.Range("J40").Value = 7666.6665
$ws.Range("K40").Value = 3333
$ws.Range("L40").Value = 7666.6665
$ws.Range("M40").Value = -3158
$ws.Range("N40").Value = -8016.6665
$ws.Range("H41").Value = 423.45456
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 423.45456
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 423.45456
$ws.Range("M41").Value = $null
$ws.Range("N41").Value = -1303.45456
$ws.Range("H48").Value = 1500
$ws.Range("I48").Value = 1500
$ws.Range("K48").Value = 4500
$ws.Range("M48").Value = -4208
$ws.Range("H56").Value = 1500
$ws.Range("I56").Value = 1500
$ws.Range("K56").Value = 4500
$ws.Range("M56").Value = -3966
$ws.Range("H64").Value = 13181.091
$ws.Range("J64").Value = 18499
$ws.Range("L64").Value = 18499
$ws.Range("N64").Value = -18995
$ws.Range("H67").Value = 13181.091
$ws.Range("J67").Value = 18499
$ws.Range("L67").Value = 18499
$ws.Range("N67").Value = -20215
$ws.Range("H116").Value = 7500
$ws.Range("I116").Value = 5000
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = -1558
$ws.Range("N116").Value = -16884
$ws.Range("H132").Value = 3966.5264
$ws.Range("I132").Value = 1112.5714
$ws.Range("J132").Value = 11957.6
$ws.Range("K132").Value = 3337.7142
$ws.Range("L132").Value = 35872.8
$ws.Range("M132").Value = -807.7142000000003
$ws.Range("N132").Value = -40932.8
$ws.Range("H135").Value = 2005
$ws.Range("I135").Value = 2005
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 18045
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -15510
$ws.Range("N135").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 500
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = $null
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H32").Value = 4552
$ws.Range("I32").Value = 4552
$ws.Range("K32").Value = 4552
$ws.Range("M32").Value = -4265
$ws.Range("H45").Value = 2918.8572
$ws.Range("I45").Value = 1486.4
$ws.Range("K45").Value = 1486.4
$ws.Range("M45").Value = -1109.4
$ws.Range("H108").Value = 105977.5
$ws.Range("J108").Value = 105977.5
$ws.Range("L108").Value = 105977.5
$ws.Range("N108").Value = -113657.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 49997
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").Value = $null
$ws.Range("H105").Value = 1684.25
$ws.Range("I105").Value = 1194.8
$ws.Range("K105").Value = 1194.8
$ws.Range("M105").Value = 552.2
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1761.375
$ws.Range("J12").Value = 2260.4
$ws.Range("L12").Value = 6781.200000000001
$ws.Range("N12").Value = -7127.200000000001
$ws.Range("H40").Value = 223
$ws.Range("I40").Value = 161.2
$ws.Range("J40").Value = 377.5
$ws.Range("K40").Value = 644.8
$ws.Range("L40").Value = 1510
$ws.Range("M40").Value = -575.8
$ws.Range("N40").Value = -1648
$ws.Range("H51").Value = 3234.5
$ws.Range("I51").Value = 901.4
$ws.Range("J51").Value = 14900
$ws.Range("K51").Value = 2704.2
$ws.Range("L51").Value = 44700
$ws.Range("M51").Value = -2244.2
$ws.Range("N51").Value = -45620
$ws.Range("H81").Value = 10000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 30000
$ws.Range("M81").Value = $null
$ws.Range("N81").Value = -32246
$ws.Range("H84").Value = 10000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 90000
$ws.Range("M84").Value = $null
$ws.Range("N84").Value = -101232
$ws.Range("H97").Value = 574.5
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 432.66666
$ws.Range("K97").Value = 3000
$ws.Range("L97").Value = 1297.99998
$ws.Range("M97").Value = -2504
$ws.Range("N97").Value = -2289.99998
$ws.Range("H122").Value = 1125.5
$ws.Range("I122").Value = 500
$ws.Range("J122").Value = 1250.6
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 11255.4
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -16155.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8779.799999999999
$ws.Range("I70").Value = 8779.799999999999
$ws.Range("K70").Value = 8779.799999999999
$ws.Range("M70").Value = -8509.799999999999
$ws.Range("H73").Value = 8779.799999999999
$ws.Range("I73").Value = 8779.799999999999
$ws.Range("K73").Value = 8779.799999999999
$ws.Range("M73").Value = -7843.799999999999
$ws.Range("H101").Value = 24499.666
$ws.Range("J101").Value = 24499.666
$ws.Range("L101").Value = 24499.666
$ws.Range("N101").Value = -30989.666
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = $null
$ws.Range("H132").Value = 2229.077
$ws.Range("I132").Value = 1633.8636
$ws.Range("J132").Value = 5502.75
$ws.Range("K132").Value = 4901.5908
$ws.Range("L132").Value = 16508.25
$ws.Range("M132").Value = -2371.5908
$ws.Range("N132").Value = -21568.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 437
$ws.Range("I22").Value = 437
$ws.Range("K22").Value = 437
$ws.Range("M22").Value = -142
$ws.Range("H27").Value = 437
$ws.Range("I27").Value = 437
$ws.Range("K27").Value = 437
$ws.Range("M27").Value = -330
$ws.Range("H94").Value = 75000
$ws.Range("J94").Value = 75000
$ws.Range("L94").Value = 75000
$ws.Range("N94").Value = -76352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 19999.5
$ws.Range("J4").Value = 19999.5
$ws.Range("L4").Value = 19999.5
$ws.Range("N4").Value = -20225.5
$ws.Range("H136").Value = 4569.6
$ws.Range("I136").Value = 2416.3333
$ws.Range("K136").Value = 7248.999899999999
$ws.Range("M136").Value = -4698.999899999999
